# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string (row 1)
$ws.Range("A1").Value = "Datos actualizados a 24 de Abril de 2020 a las 21:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 909116
$ws.Range("C4").Value = 22674
$ws.Range("D4").Value = 92266
$ws.Range("E4").Value = 765508
$ws.Range("G4").Value = 1106
$ws.Range("H4").Value = 51342

# Row 18 - Suiza
$ws.Range("E18").Value = 6088
$ws.Range("G18").Value = 40
$ws.Range("H18").Value = 1589

# Row 68 - Oman
$ws.Range("E68").Value = 1455
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 10

# Row 70 - Irak
$ws.Range("B70").Value = 1708
$ws.Range("C70").Value = 31
$ws.Range("D70").Value = 1204
$ws.Range("E70").Value = 418
$ws.Range("G70").Value = 3
$ws.Range("H70").Value = 86
